# Add a new slide (8th) to the deck by duplicating the last slide ("Amen,
# amen, amen" / Slide 7) and then clearing the two lines of sung text so the
# new slide keeps the identical placeholder, formatting (font, size, color,
# line-spacing, etc.) but is ready for new lyrics - mirroring the trailing
# pair of already-blank paragraphs that Slide 7 itself ends with.

$p = $ppt.ActivePresentation

# Slide 7 is "Amen, amen, amen" - the last slide in the deck today.
$sourceSlide = $p.Slides.Item(7)

# Duplicate() inserts the copy immediately after the source slide, so it
# naturally lands at index 8 / becomes the new last slide.
$sourceSlide.Duplicate() | Out-Null

$newSlide = $p.Slides.Item($p.Slides.Count)

# The slide has a single placeholder shape (the "Rectangle 2" subtitle
# textbox) carrying 4 paragraphs: 2 with "Amen, amen, amen" text followed by
# 2 already-empty paragraphs. Remove the two filled-in paragraphs so only the
# two blank (but fully-formatted) paragraphs remain, ready for new lyrics.
$textRange = $newSlide.Shapes.Item(1).TextFrame.TextRange
$textRange.Paragraphs(1).Delete()
$textRange.Paragraphs(1).Delete()
